$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 28
$ws.Range("I11").Value = 28
$ws.Range("K11").Value = 28
$ws.Range("M11").Value = 112
$ws.Range("H39").Value = 138.16667
$ws.Range("I39").Value = 102.57143
$ws.Range("J39").Value = 188
$ws.Range("K39").Value = 307.71429
$ws.Range("L39").Value = 564
$ws.Range("M39").Value = -11.71429000000001
$ws.Range("N39").Value = -1156
$ws.Range("H70").Value = 1125.7646
$ws.Range("I70").Value = 1067.091
$ws.Range("J70").Value = 1233.3334
$ws.Range("K70").Value = 3201.273
$ws.Range("L70").Value = 3700.0002
$ws.Range("M70").Value = -2931.273
$ws.Range("N70").Value = -4240.0002
$ws.Range("H73").Value = 1125.7646
$ws.Range("I73").Value = 1067.091
$ws.Range("J73").Value = 1233.3334
$ws.Range("K73").Value = 3201.273
$ws.Range("L73").Value = 3700.0002
$ws.Range("M73").Value = -2265.273
$ws.Range("N73").Value = -5572.0002
$ws.Range("H88").Value = 1588106.4
$ws.Range("I88").Value = 503
$ws.Range("J88").Value = 1764506.8
$ws.Range("K88").Value = 503
$ws.Range("L88").Value = 1764506.8
$ws.Range("M88").Value = -97
$ws.Range("N88").Value = -1765318.8
$ws.Range("H91").Value = 1588106.4
$ws.Range("I91").Value = 503
$ws.Range("J91").Value = 1764506.8
$ws.Range("K91").Value = 503
$ws.Range("L91").Value = 1764506.8
$ws.Range("M91").Value = 901
$ws.Range("N91").Value = -1767314.8
$ws.Range("H111").Value = 663.46155
$ws.Range("I111").Value = 759.375
$ws.Range("J111").Value = 510
$ws.Range("K111").Value = 2278.125
$ws.Range("L111").Value = 1530
$ws.Range("M111").Value = 788.875
$ws.Range("N111").Value = -7664
$ws.Range("H113").Value = 3199.5454
$ws.Range("I113").Value = 5601.6665
$ws.Range("J113").Value = 2298.75
$ws.Range("K113").Value = 5601.6665
$ws.Range("L113").Value = 2298.75
$ws.Range("M113").Value = -2347.6665
$ws.Range("N113").Value = -8806.75
$ws.Range("H125").Value = 1324.7222
$ws.Range("I125").Value = 1167.5
$ws.Range("J125").Value = 1875
$ws.Range("K125").Value = 10507.5
$ws.Range("L125").Value = 16875
$ws.Range("M125").Value = -8047.5
$ws.Range("N125").Value = -21795
$ws.Range("H132").Value = 668760.9399999999
$ws.Range("I132").Value = 2185.6155
$ws.Range("J132").Value = 5001500.5
$ws.Range("K132").Value = 6556.8465
$ws.Range("L132").Value = 15004501.5
$ws.Range("M132").Value = -4026.8465
$ws.Range("N132").Value = -15009561.5
$ws.Range("H133").Value = 43555.555
$ws.Range("J133").Value = 43555.555
$ws.Range("L133").Value = 43555.555
$ws.Range("N133").Value = -53675.555
$ws.Range("H137").Value = 2121.476
$ws.Range("I137").Value = 1770.0667
$ws.Range("J137").Value = 3000
$ws.Range("K137").Value = 5310.2001
$ws.Range("L137").Value = 9000
$ws.Range("M137").Value = -2760.2001
$ws.Range("N137").Value = -14100

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2263.3
$ws.Range("I45").Value = 1987
$ws.Range("J45").Value = 4750
$ws.Range("K45").Value = 1987
$ws.Range("L45").Value = 4750
$ws.Range("M45").Value = -1610
$ws.Range("N45").Value = -5504
$ws.Range("H61").Value = 2439.75
$ws.Range("I61").Value = 2500.5715
$ws.Range("J61").Value = 2014
$ws.Range("K61").Value = 2500.5715
$ws.Range("L61").Value = 2014
$ws.Range("M61").Value = -2288.5715
$ws.Range("N61").Value = -2438
$ws.Range("H74").Value = 1194.16
$ws.Range("I74").Value = 890.875
$ws.Range("J74").Value = 1733.3334
$ws.Range("K74").Value = 890.875
$ws.Range("L74").Value = 1733.3334
$ws.Range("M74").Value = -16.875
$ws.Range("N74").Value = -3481.3334
$ws.Range("H77").Value = 1194.16
$ws.Range("I77").Value = 890.875
$ws.Range("J77").Value = 1733.3334
$ws.Range("K77").Value = 4454.375
$ws.Range("L77").Value = 8666.666999999999
$ws.Range("M77").Value = -86.375
$ws.Range("N77").Value = -17402.667
$ws.Range("H132").Value = 2836.8276
$ws.Range("I132").Value = 2450.3809
$ws.Range("J132").Value = 3851.25
$ws.Range("K132").Value = 7351.1427
$ws.Range("L132").Value = 11553.75
$ws.Range("M132").Value = -4821.1427
$ws.Range("N132").Value = -16613.75
$ws.Range("H136").Value = 2439.75
$ws.Range("I136").Value = 2500.5715
$ws.Range("J136").Value = 2014
$ws.Range("K136").Value = 7501.7145
$ws.Range("L136").Value = 6042
$ws.Range("M136").Value = -4951.7145
$ws.Range("N136").Value = -11142

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2546.9524
$ws.Range("I86").Value = 1999
$ws.Range("J86").Value = 3642.8572
$ws.Range("K86").Value = 1999
$ws.Range("L86").Value = 3642.8572
$ws.Range("M86").Value = -876
$ws.Range("N86").Value = -5888.8572
$ws.Range("H89").Value = 2546.9524
$ws.Range("I89").Value = 1999
$ws.Range("J89").Value = 3642.8572
$ws.Range("K89").Value = 9995
$ws.Range("L89").Value = 18214.286
$ws.Range("M89").Value = -4379
$ws.Range("N89").Value = -29446.286
$ws.Range("H133").Value = 34390
$ws.Range("J133").Value = 34390
$ws.Range("L133").Value = 34390
$ws.Range("N133").Value = -44510

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 185.66667
$ws.Range("I7").Value = 185.66667
$ws.Range("K7").Value = 185.66667
$ws.Range("M7").Value = -72.66667000000001
$ws.Range("H62").Value = 9206.166999999999
$ws.Range("I62").Value = 9589.444
$ws.Range("J62").Value = 8822.888999999999
$ws.Range("K62").Value = 9589.444
$ws.Range("L62").Value = 8822.888999999999
$ws.Range("M62").Value = -8965.444
$ws.Range("N62").Value = -10070.889
$ws.Range("H65").Value = 9206.166999999999
$ws.Range("I65").Value = 9589.444
$ws.Range("J65").Value = 8822.888999999999
$ws.Range("K65").Value = 47947.22
$ws.Range("L65").Value = 44114.44499999999
$ws.Range("M65").Value = -44827.22
$ws.Range("N65").Value = -50354.44499999999
$ws.Range("H107").Value = 380.5926
$ws.Range("I107").Value = 299.16666
$ws.Range("J107").Value = 403.85715
$ws.Range("K107").Value = 299.16666
$ws.Range("L107").Value = 403.85715
$ws.Range("M107").Value = 1620.83334
$ws.Range("N107").Value = -4243.85715
$ws.Range("H132").Value = 2418.9285
$ws.Range("I132").Value = 1452
$ws.Range("J132").Value = 4159.4
$ws.Range("K132").Value = 4356
$ws.Range("L132").Value = 12478.2
$ws.Range("M132").Value = -1826
$ws.Range("N132").Value = -17538.2
$ws.Range("H134").Value = 1678.2778
$ws.Range("I134").Value = 1767.4
$ws.Range("K134").Value = 5302.200000000001
$ws.Range("M134").Value = -2767.200000000001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 32.5625
$ws.Range("I12").Value = 36.666668
$ws.Range("J12").Value = 31.615385
$ws.Range("K12").Value = 110.000004
$ws.Range("L12").Value = 94.846155
$ws.Range("M12").Value = 62.999996
$ws.Range("N12").Value = -440.846155
$ws.Range("H68").Value = 913.4286
$ws.Range("I68").Value = 849.75
$ws.Range("J68").Value = 998.3333
$ws.Range("K68").Value = 2549.25
$ws.Range("L68").Value = 2994.9999
$ws.Range("M68").Value = -1738.25
$ws.Range("N68").Value = -4616.9999
$ws.Range("H71").Value = 913.4286
$ws.Range("I71").Value = 849.75
$ws.Range("J71").Value = 998.3333
$ws.Range("K71").Value = 7647.75
$ws.Range("L71").Value = 8984.9997
$ws.Range("M71").Value = -3591.75
$ws.Range("N71").Value = -17096.9997
$ws.Range("H92").Value = 444.4
$ws.Range("I92").Value = 326
$ws.Range("J92").Value = 523.3333
$ws.Range("K92").Value = 978
$ws.Range("L92").Value = 1569.9999
$ws.Range("M92").Value = 270
$ws.Range("N92").Value = -4065.9999
$ws.Range("H109").Value = 1967.0952
$ws.Range("I109").Value = 367.66666
$ws.Range("J109").Value = 3166.6667
$ws.Range("K109").Value = 1102.99998
$ws.Range("L109").Value = 9500.000100000001
$ws.Range("M109").Value = -62.99998000000005
$ws.Range("N109").Value = -11580.0001
$ws.Range("H131").Value = 59211.39
$ws.Range("I131").Value = 168588.33
$ws.Range("J131").Value = 4522.9165
$ws.Range("K131").Value = 505764.99
$ws.Range("L131").Value = 13568.7495
$ws.Range("M131").Value = -500724.99
$ws.Range("N131").Value = -23648.7495
$ws.Range("H132").Value = 1250549.9
$ws.Range("I132").Value = 627.7143
$ws.Range("K132").Value = 5649.428699999999
$ws.Range("M132").Value = -3119.428699999999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3202.0833
$ws.Range("I80").Value = 3221.3635
$ws.Range("J80").Value = 2990
$ws.Range("K80").Value = 3221.3635
$ws.Range("L80").Value = 2990
$ws.Range("M80").Value = -2223.3635
$ws.Range("N80").Value = -4986
$ws.Range("H83").Value = 3202.0833
$ws.Range("I83").Value = 3221.3635
$ws.Range("J83").Value = 2990
$ws.Range("K83").Value = 16106.8175
$ws.Range("L83").Value = 14950
$ws.Range("M83").Value = -11114.8175
$ws.Range("N83").Value = -24934
$ws.Range("H132").Value = 2763.3333
$ws.Range("I132").Value = 1941.8235
$ws.Range("J132").Value = 4758.4287
$ws.Range("K132").Value = 5825.470499999999
$ws.Range("L132").Value = 14275.2861
$ws.Range("M132").Value = -3295.470499999999
$ws.Range("N132").Value = -19335.2861

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 442.07693
$ws.Range("I46").Value = 407
$ws.Range("J46").Value = 457.66666
$ws.Range("K46").Value = 407
$ws.Range("L46").Value = 457.66666
$ws.Range("M46").Value = -219
$ws.Range("N46").Value = -833.66666
$ws.Range("H132").Value = 2209.5715
$ws.Range("I132").Value = 1856.88
$ws.Range("J132").Value = 3091.3
$ws.Range("K132").Value = 5570.64
$ws.Range("L132").Value = 9273.900000000001
$ws.Range("M132").Value = -3040.64
$ws.Range("N132").Value = -14333.9

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1237.9459
$ws.Range("I132").Value = 781.1613
$ws.Range("J132").Value = 3598
$ws.Range("K132").Value = 2343.4839
$ws.Range("L132").Value = 10794
$ws.Range("M132").Value = 186.5160999999998
$ws.Range("N132").Value = -15854

Write-Output "Applied all profit/price updates across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets"